# Actualiza la lista de lavadas: se agrego la placa en el ticket de
# salida y ahora los cajeros pueden ver la lista completa.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fila 2: actualizar datos existentes ---
$ws.Range("B2").Value = "JKJ3865"
$ws.Range("C2").Value = "NISSAN"
$ws.Range("D2").Value = "2024-06-13 13:00 PM"
$ws.Range("E2").Value = "2024-06-18 23:00:00"

# --- Fila 3: se quita el nombre de cliente "edicion" y se actualizan datos ---
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "JPK2456"
$ws.Range("C3").Value = "HONDA"
$ws.Range("D3").Value = "2024-06-13 13:09 PM"
$ws.Range("E3").Value = "2024-06-18 23:30:00"

# --- Filas nuevas 4-11: nuevos vehiculos en la lista ---
$newRows = @(
    @{ Row = 4;  A = "";                          B = "JLA3831"; C = "HONDA";      D = "2024-06-13 14:07 PM"; E = "2024-06-11 22:00:00" },
    @{ Row = 5;  A = "";                          B = "294XUJ";  C = "VOLKSWAGEN"; D = "2024-06-13 16:13 PM"; E = "2024-06-17 18:30:00" },
    @{ Row = 6;  A = "";                          B = "G19BMM";  C = "MAZDA";      D = "2024-06-13 18:57 PM"; E = "" },
    @{ Row = 7;  A = "";                          B = "JSJ2237"; C = "MINI";       D = "2024-06-13 21:21 PM"; E = "" },
    @{ Row = 8;  A = "";                          B = "JNR1775"; C = "NISSAN";     D = "2024-06-14 02:45 AM"; E = "2024-06-21 19:00:00" },
    @{ Row = 9;  A = "";                          B = "JUB9526"; C = "NISSAN";     D = "2024-06-14 03:04 AM"; E = "2024-06-18 23:00:00" },
    @{ Row = 10; A = "ALEJANDRA VANESSA JIMENEZ"; B = "JLY1080"; C = "MITSUBISHI"; D = "2024-06-14 03:13 AM"; E = "" },
    @{ Row = 11; A = "";                          B = "HWF430A"; C = "PEUGEOT";    D = "2024-06-14 03:58 AM"; E = "" }
)

foreach ($r in $newRows) {
    if ($r.A -ne "") {
        $ws.Cells.Item($r.Row, 1).Value = $r.A
    }
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    if ($r.E -ne "") {
        $ws.Cells.Item($r.Row, 5).Value = $r.E
    }
}
